$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# B2 (shortname): numeric 1 -> text "t1"
$ws1.Range("B2").Value = "t1"

# B8 (principaldefault): 30000000 -> 30000
$ws1.Range("B8").Value = 30000

# B10 (nominalinterestratedefault): 14 -> 12
$ws1.Range("B10").Value = 12

# B31 (maximumallowedoutstandingbalancefortranchloan): 30000000 -> 1000000
$ws1.Range("B31").Value = 1000000

# Update the active cell selection on each sheet to reflect the diff
$ws2.Range("B20").Select()
$ws1.Range("B29").Select()
